# bars-stacked.xlsx: add a "style" / "default" metadata row to the "meta"
# sheet, pushing the previously-trailing blank (formatted) row down by one.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("meta")

# Row 5 ("A5") currently holds the trailing blank cell (style copied from the
# other key cells in column A, e.g. A4). Copy that formatting down to the new
# blank row 6 first, then turn row 5 into the new "style" / "default" entry.
$meta.Range("A4").Copy() | Out-Null
$meta.Range("A6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$meta.Range("A4").Copy() | Out-Null
$meta.Range("A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$meta.Range("A5").Value = "style"
$meta.Range("B5").Value = "default"

$excel.CutCopyMode = 0
